$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 (pushes the old row 7 and everything below it down by one)
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the "intended use" label/value pair
$ws.Range("A7").Value = "intendedUse"
$ws.Range("B7").Value = "Epi-validated outbreak"

# Seed A7/B7 formatting from the header row (row 9 after the insert) so the
# new cells pick up applyFont/applyBorder flags consistent with the rest of
# the sheet, then tweak the font to the bold/plain 12pt black Calibri look.
$ws.Range("A9").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("B9").Copy()
$ws.Range("B7").PasteSpecial(-4122)

$ws.Range("A7").Font.Color = 0
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Font.Size = 12

$ws.Range("B7").Font.Color = 0
$ws.Range("B7").Font.Bold = $false
$ws.Range("B7").Font.Size = 12
$ws.Range("B7").HorizontalAlignment = -4131

$wb.Application.CutCopyMode = $false

# Update the sheet view/selection to match the edited state
$ws.Range("A7:B7").Select()
